$wb = $excel.ActiveWorkbook

# Helper: write a text value into a cell while avoiding Excel's automatic
# "looks like a number" -> numeric type coercion (which would otherwise
# turn these string cells into real numeric cells and also permanently
# tag them with a Text-format style). We briefly force Text format, set
# the value, then reset the cell style back to Normal/General so the
# final cell carries no explicit style index (matching the source file,
# where every data cell is unstyled).
function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------
# Sheet "Restricciones_del_follower" (3rd sheet)
# Columns: A=Expression  B=Function_Evaluation  C=Restriction_Set_Type
#          D=Lambda_value  E=Beta_value  F=Gamma_value
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Restricciones_del_follower")

Set-TextValue $ws3.Range("A2") "1.1000000000000056 - 2x_1 + y_1 - y_2"
Set-TextValue $ws3.Range("B2") "1.3999999999999944"
Set-TextValue $ws3.Range("D2") "0.92"
Set-TextValue $ws3.Range("F2") "0"

Set-TextValue $ws3.Range("A3") "2.5499999999999927 + x_1 - 3x_2 + y_2"
Set-TextValue $ws3.Range("B3") "-4.549999999999993"
Set-TextValue $ws3.Range("D3") "0.36"
Set-TextValue $ws3.Range("E3") "0"

Set-TextValue $ws3.Range("A4") "102.2 - y_1"
Set-TextValue $ws3.Range("B4") "-102.2"
Set-TextValue $ws3.Range("D4") "0.49"
Set-TextValue $ws3.Range("E4") "0"
Set-TextValue $ws3.Range("F4") "0"

Set-TextValue $ws3.Range("A5") "-0.2 - y_2"
Set-TextValue $ws3.Range("B5") "-0.2"
Set-TextValue $ws3.Range("D5") "0.75"
Set-TextValue $ws3.Range("E5") "-2.4"
Set-TextValue $ws3.Range("F5") "-5.300000000000001"

# ---------------------------------------------------------------
# Sheet "Punto_modificado" (4th sheet)
# Columns: A=x_1  B=x_2  C=y_1  D=y_2
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Punto_modificado")

Set-TextValue $ws4.Range("A2") "51.550000000000004"
Set-TextValue $ws4.Range("B2") "18.099999999999998"
Set-TextValue $ws4.Range("C2") "102.2"
Set-TextValue $ws4.Range("D2") "0.2"

# ---------------------------------------------------------------
# Sheet "Vector_bf" (5th sheet)
# ---------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Vector_bf")

Set-TextValue $ws5.Range("A2") "3.57"
Set-TextValue $ws5.Range("A3") "0.31000000000000005"

# ---------------------------------------------------------------
# Sheet "Vector_BF" (6th sheet)
# ---------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Vector_BF")

Set-TextValue $ws6.Range("A2") "2.0"
Set-TextValue $ws6.Range("A3") "-1.0"
Set-TextValue $ws6.Range("A4") "-0.5"
Set-TextValue $ws6.Range("A5") "-2.4"

# "Vector_Alpha" (7th sheet) and all other sheets/cells are unchanged.
